$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New data rows 43-51 appended below the existing table (which ended at 42).
# Column layout: A = 샘플명 (text or, for row 50, a plain number code),
# B..K = the 10 numeric measurement columns.
# Odd "group" rows reuse the banded format of row 41 (style ids 6/7),
# even "group" rows reuse the banded format of row 42 (style ids 8/9).
# ---------------------------------------------------------------------------

function Set-RowValues($row, $a, $values) {
    $ws.Range("A$row").Value = $a
    $cols = @("B","C","D","E","F","G","H","I","J","K")
    for ($i = 0; $i -lt $values.Count; $i++) {
        $ws.Cells.Item($row, $i + 2).Value = $values[$i]
    }
}

# Row 43 (odd banding -> copy format from row 41)
$ws.Range("A41:K41").Copy()
$ws.Range("A43").PasteSpecial(-4122)
$excel.CutCopyMode = $false
Set-RowValues 43 "GRAY(K25B003)" @(64.3, -1.41, -3.33, 92.09, 4.489, 0.01, 3.37, 0, 0, 0.04)

# Row 44 (even banding -> copy format from row 42)
$ws.Range("A42:K42").Copy()
$ws.Range("A44").PasteSpecial(-4122)
$excel.CutCopyMode = $false
Set-RowValues 44 "RAL5005" @(33.18, 2.75, -28.2, 17.73, 8.99, 0, 0.24, 67.24, 0, 5.81)

# Row 45 (odd banding)
$ws.Range("A41:K41").Copy()
$ws.Range("A45").PasteSpecial(-4122)
$excel.CutCopyMode = $false
Set-RowValues 45 "3.3B 4.1/4.4" @(42.2, -12.28, -14.95, 37.59, 8.23, 0, 6.5, 38.54, 21.1, 0)

# Row 46 (even banding)
$ws.Range("A42:K42").Copy()
$ws.Range("A46").PasteSpecial(-4122)
$excel.CutCopyMode = $false
Set-RowValues 46 "tn726" @(52.88, -7.3, -9.64, 74.56, 5.84, 3.33, 4.9, 11.37, 0, 0)

# Row 47 (odd banding)
$ws.Range("A41:K41").Copy()
$ws.Range("A47").PasteSpecial(-4122)
$excel.CutCopyMode = $false
Set-RowValues 47 "6.25yr" @(60.4, 11.53, 19.3, 54.6, 46.47, 6.56, 1.37, 0, 0, 0)

# Row 48 (even banding)
$ws.Range("A42:K42").Copy()
$ws.Range("A48").PasteSpecial(-4122)
$excel.CutCopyMode = $false
Set-RowValues 48 "2.5pb" @(50.81, 1.4, 1.98, 95.89, 4.06, 0.01, 0.03, 0, 0, 0)

# Row 49 (odd banding)
$ws.Range("A41:K41").Copy()
$ws.Range("A49").PasteSpecial(-4122)
$excel.CutCopyMode = $false
Set-RowValues 49 "WHITE" @(88.65, -0.88, 5.33, 95.75, 4.21, 0.03, 0.01, 0, 0, 0)

# Row 50 (even banding, but A50 is left-aligned and holds a plain number)
$ws.Range("A42:K42").Copy()
$ws.Range("A50").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A50").HorizontalAlignment = -4131
$ws.Range("A50").Value = 3041
$cols50 = @(88.13, -0.97, 3.44, 98.48, 1.38, 0.06, 0.08, 0, 0, 0)
for ($i = 0; $i -lt $cols50.Count; $i++) {
    $ws.Cells.Item(50, $i + 2).Value = $cols50[$i]
}

# Row 51 (odd banding for B:K, but A51 drops back to the default/no style)
$ws.Range("A41:K41").Copy()
$ws.Range("A51").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A51").ClearFormats()
Set-RowValues 51 "WHITE" @(88.38, 1.29, 4.27, 97.93, 1.71, 0.29, 0.06, 0, 0, 0)

# ---------------------------------------------------------------------------
# View state: keep the header frozen, scroll the table so the new rows are
# visible, and leave the selection where Excel would after entering this data.
# ---------------------------------------------------------------------------
[void]$ws.Range("K51").Select()
